# Re-sort the CHE grid rows (SetName / pset_co pairs) on the geo_sets sheet
# so that the numeric suffix of "rez_CHE_<n>" is in ascending numeric order
# (0,1,2,3,...,9,10,11,...) instead of the previous text/lexicographic order
# (0,1,10,11,...,15,17,18,19,2,20,...,25,3,4,...,9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 160
$endRow = 184

# Read the current (lexicographically ordered) values from columns B and C
$data = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $key = [int]([regex]::Match([string]$b, '\d+$').Value)
    $data += [PSCustomObject]@{ B = $b; C = $c; Key = $key }
}

# Sort ascending by the numeric suffix
$sorted = $data | Sort-Object -Property Key

# Write the sorted values back into the same B:C range
$i = 0
for ($r = $startRow; $r -le $endRow; $r++) {
    $row = $sorted[$i]
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $i++
}
